# Updated cryptos list on Sat Jul 13 16:52:28 UTC 2024 with GitHub Actions
#
# Refreshes prices / 1h volume-change figures from the latest
# coinranking.com snapshot, and re-ranks rows 24/25 (Polygon now ranks
# above Litecoin).
#
# NOTE: column D holds price text like "58.716.23" / "13.10" / "0.0000172".
# Excel's COM layer auto-detects plain decimal-looking strings and would
# silently coerce them to numbers (dropping significant trailing zeros,
# switching to scientific notation, etc.), so for any new value that looks
# like a plain number we force the cell to Text format first to keep it a
# literal string, matching the original inlineStr cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "58.716.23"
$ws.Range("E2").Value = "  +1.23%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "3.154.74"
$ws.Range("E3").Value = "  +0.92%  "

# --- Row 5: BNB ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.12"
$ws.Range("E5").Value = "  -0.25%  "

# --- Row 6: Solana ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.65"
$ws.Range("E6").Value = "  +1.23%  "

# --- Row 7: USDC ---
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.12%  "

# --- Row 8: XRP ---
$ws.Range("E8").Value = "  +14.52%  "

# --- Row 9: Toncoin ---
$ws.Range("E9").Value = "  +0.38%  "

# --- Row 10: Cardano ---
$ws.Range("E10").Value = "  +5.55%  "

# --- Row 11: Dogecoin ---
$ws.Range("E11").Value = "  +2.92%  "

# --- Row 12: TRON ---
$ws.Range("E12").Value = "  +2.67%  "

# --- Row 13: Wrapped liquid staked Ether 2.0 ---
$ws.Range("D13").Value = "3.697.87"
$ws.Range("E13").Value = "  +0.91%  "

# --- Row 14: Avalanche ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.04"
$ws.Range("E14").Value = "  +1.67%  "

# --- Row 15: Shiba Inu ---
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  +5.16%  "

# --- Row 16: Wrapped BTC ---
$ws.Range("D16").Value = "58.754.23"
$ws.Range("E16").Value = "  +1.15%  "

# --- Row 17: Polkadot ---
$ws.Range("E17").Value = "  +4.20%  "

# --- Row 18: Wrapped Ether ---
$ws.Range("D18").Value = "3.155.09"
$ws.Range("E18").Value = "  +0.79%  "

# --- Row 19 ---
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.10"
$ws.Range("E19").Value = "  +3.70%  "

# --- Row 20 ---
$ws.Range("E20").Value = "  +0.71%  "

# --- Row 21 ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.92"
$ws.Range("E21").Value = "  +5.37%  "

# --- Row 22 ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.80"
$ws.Range("E22").Value = "  +1.70%  "

# --- Row 23 ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.77%  "

# --- Rows 24/25: Polygon and Litecoin swap rank positions ---
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.522"
$ws.Range("E24").Value = "  +3.70%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.87"
$ws.Range("E25").Value = "  +1.15%  "

# --- Row 26 ---
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +0.19%  "

# --- Row 27 ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.21%  "

# --- Row 28 ---
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.27"
$ws.Range("E28").Value = "  +13.78%  "

# --- Row 29 ---
$ws.Range("D29").Value = "0.0₃0865"
$ws.Range("E29").Value = "  -1.03%  "

# --- Row 30 ---
$ws.Range("E30").Value = "  +0.54%  "

# --- Row 31 ---
$ws.Range("E31").Value = "  -1.12%  "

# --- Row 32 ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.10"
$ws.Range("E32").Value = "  +3.34%  "

# --- Row 33 ---
$ws.Range("E33").Value = "  +3.93%  "

# --- Row 34 ---
$ws.Range("E34").Value = "  +2.10%  "

# --- Row 35: Monero ---
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.18"
$ws.Range("E35").Value = "  +0.28%  "

# --- Row 36: Aptos ---
$ws.Range("E36").Value = "  +3.26%  "

# --- Row 37: ImmutableX ---
$ws.Range("E37").Value = "  +6.13%  "

# --- Row 38: EnergySwap ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.18"
$ws.Range("E38").Value = "  -2.35%  "

# --- Row 39: Stacks ---
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.68"
$ws.Range("E39").Value = "  -0.79%  "

# --- Row 40: Maker ---
$ws.Range("D40").Value = "2.644.94"
$ws.Range("E40").Value = "  +7.76%  "

# --- Row 41: Hedera ---
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0687"
$ws.Range("E41").Value = "  +2.88%  "

# --- Row 42: Filecoin ---
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.26"
$ws.Range("E42").Value = "  +6.94%  "

# --- Row 43: OKB ---
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.87"
$ws.Range("E43").Value = "  +3.52%  "

# --- Row 44: Mantle ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.709"
$ws.Range("E44").Value = "  +1.40%  "

# --- Row 45: VeChain ---
$ws.Range("E45").Value = "  +6.34%  "

# --- Row 46: FirstDigitalUSD ---
$ws.Range("E46").Value = "  -0.01%  "

# --- Row 47: RenzoRestakedETH ---
$ws.Range("D47").Value = "3.196.61"
$ws.Range("E47").Value = "  +0.88%  "

# --- Row 48: Stellar ---
$ws.Range("E48").Value = "  +13.79%  "

# --- Row 49: Cosmos ---
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.22"
$ws.Range("E49").Value = "  +3.13%  "

# --- Row 50: ONDO ---
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.983"
$ws.Range("E50").Value = "  +0.52%  "

# --- Row 51: InjectiveProtocol ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.41"
$ws.Range("E51").Value = "  +3.08%  "
